# ----------------------------------------------------------------------
# "starting on complete panel2 mape, skip NA when rename nodes"
#
# - panel1: just move the active selection (G12), no data changes.
# - panel2: replace the two starter rows with the full gating-tree node
#   list (35 nodes), using "NA" as the not-yet-mapped placeholder in
#   columns B/C (skipped during renaming), except for the two nodes that
#   already had a real mapping (PE-A -> "Live cells (PE-)" and
#   Singlets -> "Single Cells (FSC-H v FSC-W)").
# - widen column A on panel2 to fit the longer node names.
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- panel1: only the selection/cursor moved -------------------------
$ws1 = $wb.Worksheets.Item("panel1")
$ws1.Range("G12").Select()

# --- panel2: rebuild the node map -------------------------------------
$ws2 = $wb.Worksheets.Item("panel2")
$ws2.Select()

# Full ordered list of node names for column A, rows 2..36.
$nodeNames = @(
  "boundary",
  "CD45+",
  "PE-A+",
  "PE-A",
  "FSC-H+",
  "SingletsH",
  "FSC-W+",
  "SingletsW",
  "Singlets",
  "nonDebris",
  "PBMC",
  "CD19+",
  "CD19-",
  "CD3+",
  "CD3-",
  "D_NK_M",
  "CD14+",
  "CD14+/CD16+",
  "CD16-",
  "CD14-",
  "CD20-",
  "HLA-DR+",
  "Dendritic",
  "BB515-A-BV 711-A+",
  "BB515-A+BV 711-A+",
  "BB515-A+BV 711-A-",
  "BB515-A-BV 711-A-",
  "CD56+",
  "CD20-/CD16+",
  "CD16-CD56-",
  "CD16+CD56-",
  "CD16-CD56+",
  "CD16+CD56+",
  "CD56PlusPlus",
  "CD45-"
)

# Column C overrides: node name -> already-known mapped description.
# Everything else in B/C gets the "NA" skip placeholder.
$mappedC = @{
  "PE-A"     = "Live cells (PE-)"
  "Singlets" = "Single Cells (FSC-H v FSC-W)"
}

# Wipe any leftover formatting on the old starter rows (2:3) and on the
# old blank placeholder cells (C20:C27) so the rebuilt rows start clean.
$ws2.Range("A2:C4").ClearFormats()
$ws2.Range("C20:C27").ClearFormats()

# Give column C (rows 5 and down) the same "Manual" style as the header
# (C1), matching the look of the original two mapped rows.
$ws2.Range("C1").Copy()
$ws2.Range("C5:C36").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column A first, top to bottom, for every row ...
for ($i = 0; $i -lt $nodeNames.Length; $i++) {
  $r = $i + 2
  $ws2.Cells.Item($r, 1).Value = $nodeNames[$i]
}

# ... then columns B/C, so the "NA" skip placeholder is the very last
# new shared string introduced (matches the authoring order: name the
# nodes, then go back and mark everything not yet mapped as NA).
for ($i = 0; $i -lt $nodeNames.Length; $i++) {
  $r = $i + 2
  $name = $nodeNames[$i]

  $ws2.Cells.Item($r, 2).Value = "NA"

  if ($mappedC.ContainsKey($name)) {
    $ws2.Cells.Item($r, 3).Value = $mappedC[$name]
  } else {
    $ws2.Cells.Item($r, 3).Value = "NA"
  }
}

# Widen column A to fit the longest node name.
$ws2.Columns.Item(1).ColumnWidth = 16.666666666666668

$ws2.Range("A12").Select()
